$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 1063, shifting existing rows 1063-1153 down to 1064-1154
$ws.Rows.Item(1063).Insert()

# Fill in the new row 1063 with the new weekly record
$ws.Cells.Item(1063, 1).Value = 3
$ws.Cells.Item(1063, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1063, 3).Value = "Coquimbo"
$ws.Cells.Item(1063, 4).Value = 45166
$ws.Cells.Item(1063, 5).Value = 5
$ws.Cells.Item(1063, 6).Value = 100112024
$ws.Cells.Item(1063, 7).Value = "Choclo"
$ws.Cells.Item(1063, 8).Value = "Dulce o Americano"
$ws.Cells.Item(1063, 9).Value = "Primera"
$ws.Cells.Item(1063, 10).Value = 60
$ws.Cells.Item(1063, 11).Value = 45000
$ws.Cells.Item(1063, 12).Value = 45000
$ws.Cells.Item(1063, 13).Value = 45000
$ws.Cells.Item(1063, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(1063, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1063, 16).Value = 643
$ws.Cells.Item(1063, 17).Value = 70
$ws.Cells.Item(1063, 18).Value = "Hortaliza"
